$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @("2:6.00","3:1.00","4:1.00","5:1.00","6:1.00","7:1.00","8:1.00","9:1.00","10:1.00","11:2.00","12:1.00","13:1.00","14:1.00","15:8.00","16:1.00","17:1.00","18:6.00","19:4.00","20:3.00","21:4.00","22:1.00","23:1.00","24:14.00","25:1.00","26:1.00","27:1.00","28:14.00","29:2.00","30:1.00","31:3.00","32:1.00","33:1.00","34:1.00","35:1.00","36:4.00","37:9.00","38:1.00","39:1.00","40:1.00","41:1.00","42:1.00","43:2.00","44:1.00","45:2.00","46:4.00","47:1.00","48:1.00","49:1.00","50:2.00","51:2.00","52:4.00","53:6.00","54:10.00","55:1.00","56:1.00","57:4.00","58:14.00","59:2.00","60:5.00","61:1.00","62:2.00","63:1.00","64:1.00","65:1.00","66:2.00","67:1.00","68:1.00","69:1.00","70:2.00","71:1.00","72:1.00","73:2.00","74:2.00","75:2.00","76:1.00","77:1.00","78:1.00","79:1.00","80:1.00","81:1.00","82:1.00","83:1.00","84:1.00","85:12.00","86:1.00","87:5.00","88:150.00","89:1.00","90:4.00","91:2.00","92:1.00","93:1.00","94:14.00","95:4.00","96:1.00","97:8.00","98:1.00","99:3.00","100:1.00","101:110.00","102:33.00","103:36.00","104:2.00","105:1.00","106:2.00","107:5.00","108:1.00","109:1.00","110:200.00","111:10.00","112:2.00","113:1.00","114:1.00","115:1.00","116:1.00","117:1.00","118:2.00","119:2.00","120:1.00","121:1.00","122:2.00","123:1.00","124:1.00","125:2.00","126:8.00","127:1.00","128:2.00","129:1.00","130:1.00","131:2.00","132:1.00","133:2.00","134:3.00","135:1.00","136:1.00","137:7.00","138:1.00","139:2.00","140:1.00","141:1.00","142:4.00","143:1.00","144:1.00","145:1.00","146:4.00","147:3.00","148:1.00","149:1.00","150:5.00","151:1.00","152:1.00","153:1.00","154:1.00","155:1.00","156:1.00","157:1.00","158:1.00","159:1.00","160:1.00","161:1.00","162:1.00","163:70.00","164:140.00","165:70.00","166:70.00","167:1.00","168:1.00","169:1.00","170:1.00","171:4.00","172:5.00","173:4.00","174:6.00","175:4.00","176:1.00","177:4.00","178:4.00","179:16.00","180:4.00","181:1.00","182:5.00","183:5.00","184:2.00","185:3.00","186:4.00","187:6.00","188:4.00","189:5.00","190:3.00","191:1.00","192:1.00","193:6.00","194:1.00","195:1.00","196:2.00","197:1.00","198:1.00","199:1.00","200:1.00","201:1.00","202:1.00","203:1.00","204:35.00","205:50.00","206:13.00","207:4.00","208:22.00","209:25.00","210:33.00","211:10.00","212:18.00","213:22.00","214:6.00","215:16.00","216:28.00","217:175.00","218:1.00","219:1.00","220:1.00","221:1.00","222:1.00","223:1.00","224:1.00","225:1.00","226:155.00","227:20.00","228:13.00","229:8.00","230:4.00","231:1.00","232:41.00","233:18.00","234:10.00","235:6.00","236:3.00","237:1.00","238:1.00","239:1.00","240:1.00","241:1.00","242:1.00","243:1.00","244:1.00","245:1.00","246:1.00","247:1.00","248:1.00","249:1.00","250:1.00","251:1.00","252:1.00","253:1.00","254:1.00","255:6.00","256:1.00","257:1.00","258:8.00","259:14.00","260:9.00","261:1.00","262:1.00","263:1.00","264:1.00","265:1.00","266:1.00","267:1.00","268:1.00","269:1.00","270:1.00","271:1.00","272:1.00","273:1.00","274:1.00","275:50.00","276:26.00","277:3.00","278:1.00","279:1.00","280:4.00","281:2.00","282:12.00","283:1.00","284:1.00","285:6.00","286:4.00","287:1.00","288:1.00","289:50.00","290:1.00","291:8.00","292:36.00","293:1.00","294:1.00","295:50.00","296:1.00","297:8.00","298:4.00","299:3.00","300:4.00","301:1.00","302:1.00")

foreach ($item in $data) {
    $parts = $item.Split(":")
    $row = [int]$parts[0]
    $val = $parts[1]
    $cell = $ws.Cells.Item($row, 13)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}
